# Applies the Titan_Profits sheet data refresh (currentAveragePrice* and derived
# Leve price/profit columns) produced by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3153.5715
$ws.Range("J17").Value = 3153.5715
$ws.Range("L17").Value = 9460.7145
$ws.Range("N17").Value = -9796.7145

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 956.3333
$ws.Range("I39").Value = 1462.125
$ws.Range("J39").Value = 378.2857
$ws.Range("K39").Value = 4386.375
$ws.Range("L39").Value = 1134.8571
$ws.Range("M39").Value = -4090.375
$ws.Range("N39").Value = -1726.8571

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2377413.2
$ws.Range("I62").Value = 3377498
$ws.Range("K62").Value = 3377498
$ws.Range("M62").Value = -3376874

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2377413.2
$ws.Range("I65").Value = 3377498
$ws.Range("K65").Value = 16887490
$ws.Range("M65").Value = -16884370

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M97").ClearContents()
$ws.Range("H97").Value = 3750
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 3750
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 11250
$ws.Range("N97").Value = -12242

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4945.4287
$ws.Range("I113").Value = 3168.3333
$ws.Range("K113").Value = 3168.3333
$ws.Range("M113").Value = 85.66670000000022

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1095.9756
$ws.Range("J129").Value = 1159.5526
$ws.Range("L129").Value = 3478.6578
$ws.Range("N129").Value = -13478.6578

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6779565
$ws.Range("I138").Value = 1896116.5
$ws.Range("J138").Value = 10002641
$ws.Range("K138").Value = 5688349.5
$ws.Range("L138").Value = 30007923
$ws.Range("M138").Value = -5683209.5
$ws.Range("N138").Value = -30018203

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2244
$ws.Range("I122").Value = 2170.4614
$ws.Range("J122").Value = 3200
$ws.Range("K122").Value = 6511.3842
$ws.Range("L122").Value = 9600
$ws.Range("M122").Value = -4061.3842
$ws.Range("N122").Value = -14500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 1101.3334
$ws.Range("I8").Value = 652
$ws.Range("J8").Value = 2000
$ws.Range("K8").Value = 652
$ws.Range("L8").Value = 2000
$ws.Range("M8").Value = -512
$ws.Range("N8").Value = -2280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1364.6111
$ws.Range("I20").Value = 1229.5
$ws.Range("J20").Value = 1634.8334
$ws.Range("K20").Value = 1229.5
$ws.Range("L20").Value = 1634.8334
$ws.Range("M20").Value = -982.5
$ws.Range("N20").Value = -2128.8334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1034.7059
$ws.Range("I99").Value = 756.4286
$ws.Range("J99").Value = 2333.3333
$ws.Range("K99").Value = 756.4286
$ws.Range("L99").Value = 2333.3333
$ws.Range("M99").Value = 741.5714
$ws.Range("N99").Value = -5329.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3079.9773
$ws.Range("I105").Value = 2793.6
$ws.Range("J105").Value = 3693.6428
$ws.Range("K105").Value = 2793.6
$ws.Range("L105").Value = 3693.6428
$ws.Range("M105").Value = -1046.6
$ws.Range("N105").Value = -7187.6428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5159.0786
$ws.Range("I31").Value = 1553.8334
$ws.Range("J31").Value = 10309.429
$ws.Range("K31").Value = 1553.8334
$ws.Range("L31").Value = 10309.429
$ws.Range("M31").Value = -1258.8334
$ws.Range("N31").Value = -10899.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5159.0786
$ws.Range("I34").Value = 1553.8334
$ws.Range("J34").Value = 10309.429
$ws.Range("K34").Value = 1553.8334
$ws.Range("L34").Value = 10309.429
$ws.Range("M34").Value = -1351.8334
$ws.Range("N34").Value = -10713.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6404.95
$ws.Range("I70").Value = 6382.294
$ws.Range("J70").Value = 6533.3335
$ws.Range("K70").Value = 6382.294
$ws.Range("L70").Value = 6533.3335
$ws.Range("M70").Value = -6112.294
$ws.Range("N70").Value = -7073.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6404.95
$ws.Range("I73").Value = 6382.294
$ws.Range("J73").Value = 6533.3335
$ws.Range("K73").Value = 6382.294
$ws.Range("L73").Value = 6533.3335
$ws.Range("M73").Value = -5446.294
$ws.Range("N73").Value = -8405.333500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1729.2727
$ws.Range("I122").Value = 1546.2858
$ws.Range("J122").Value = 2049.5
$ws.Range("K122").Value = 4638.857400000001
$ws.Range("L122").Value = 6148.5
$ws.Range("M122").Value = -2188.857400000001
$ws.Range("N122").Value = -11048.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2632.2559
$ws.Range("I126").Value = 2288.1667
$ws.Range("K126").Value = 6864.500100000001
$ws.Range("M126").Value = -4394.500100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 43869.875
$ws.Range("J139").Value = 43869.875
$ws.Range("L139").Value = 43869.875
$ws.Range("N139").Value = -54149.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2512.361
$ws.Range("I7").Value = 1882
$ws.Range("J7").Value = 2789.72
$ws.Range("K7").Value = 1882
$ws.Range("L7").Value = 2789.72
$ws.Range("M7").Value = -1770
$ws.Range("N7").Value = -3013.72

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 709.0833
$ws.Range("I46").Value = 590
$ws.Range("J46").Value = 1066.3334
$ws.Range("K46").Value = 590
$ws.Range("L46").Value = 1066.3334
$ws.Range("M46").Value = -402
$ws.Range("N46").Value = -1442.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3539.7307
$ws.Range("I122").Value = 2614.6667
$ws.Range("J122").Value = 3817.25
$ws.Range("K122").Value = 7844.000100000001
$ws.Range("L122").Value = 11451.75
$ws.Range("M122").Value = -5394.000100000001
$ws.Range("N122").Value = -16351.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2512.361
$ws.Range("I126").Value = 1882
$ws.Range("J126").Value = 2789.72
$ws.Range("K126").Value = 5646
$ws.Range("L126").Value = 8369.16
$ws.Range("M126").Value = -3176
$ws.Range("N126").Value = -13309.16

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1412.45
$ws.Range("I122").Value = 1269.9333
$ws.Range("J122").Value = 1840
$ws.Range("K122").Value = 3809.7999
$ws.Range("L122").Value = 5520
$ws.Range("M122").Value = -1359.7999
$ws.Range("N122").Value = -10420

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 39709.92
$ws.Range("I126").Value = 67498.92999999999
$ws.Range("J126").Value = 1815.8182
$ws.Range("K126").Value = 202496.79
$ws.Range("L126").Value = 5447.4546
$ws.Range("M126").Value = -200026.79
$ws.Range("N126").Value = -10387.4546
